# Auto-generated edit script applying cryptos.xlsx price/volume/coin updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.225.06"
$ws.Range("D2").Style = "Normal"

$ws.Range("E3").Value = "  -0.38%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.272.59"
$ws.Range("D3").Style = "Normal"

$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"

$ws.Range("E5").Value = "  +0.40%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.23"
$ws.Range("D5").Style = "Normal"

$ws.Range("E6").Value = "  +2.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.04"
$ws.Range("D6").Style = "Normal"

$ws.Range("E7").Value = "  -1.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.526"
$ws.Range("D7").Style = "Normal"

$ws.Range("E8").Value = "  -0.06%  "

$ws.Range("E9").Value = "  +0.79%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.496"
$ws.Range("D9").Style = "Normal"

$ws.Range("E10").Value = "  +3.93%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.56"
$ws.Range("D10").Style = "Normal"

$ws.Range("E11").Value = "  -1.61%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0792"
$ws.Range("D11").Style = "Normal"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.113"
$ws.Range("D12").Style = "Normal"

$ws.Range("E13").Value = "  +2.20%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.86"
$ws.Range("D13").Style = "Normal"

$ws.Range("E14").Value = "  -0.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.624.27"
$ws.Range("D14").Style = "Normal"

$ws.Range("E15").Value = "  +2.80%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.80"
$ws.Range("D15").Style = "Normal"

$ws.Range("E16").Value = "  -0.89%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.259.31"
$ws.Range("D16").Style = "Normal"

$ws.Range("E17").Value = "  +0.85%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.797"
$ws.Range("D17").Style = "Normal"

$ws.Range("E18").Value = "  -0.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.088.50"
$ws.Range("D18").Style = "Normal"

$ws.Range("E19").Value = "  -3.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.47"
$ws.Range("D19").Style = "Normal"

$ws.Range("E20").Value = "  -1.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0908"
$ws.Range("D20").Style = "Normal"

$ws.Range("E21").Value = "  +0.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.03"
$ws.Range("D21").Style = "Normal"

$ws.Range("E22").Value = "  +0.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.47"
$ws.Range("D22").Style = "Normal"

$ws.Range("E23").Value = "  -2.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "238.70"
$ws.Range("D23").Style = "Normal"

$ws.Range("E24").Value = "  -1.91%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.56"
$ws.Range("D24").Style = "Normal"

$ws.Range("B25").Value = "ImmutableX"
$ws.Range("C25").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("E25").Value = "  -0.77%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.93"
$ws.Range("D25").Style = "Normal"

$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"

$ws.Range("E27").Value = "  -1.75%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.72"
$ws.Range("D27").Style = "Normal"

$ws.Range("E28").Value = "  +4.80%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "37.46"
$ws.Range("D28").Style = "Normal"

$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("E29").Value = "  -1.90%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.52"
$ws.Range("D29").Style = "Normal"

$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("E30").Value = "  +1.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.12"
$ws.Range("D30").Style = "Normal"

$ws.Range("E31").Value = "  +0.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "160.73"
$ws.Range("D31").Style = "Normal"

$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("E32").Value = "  -2.38%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.24"
$ws.Range("D32").Style = "Normal"

$ws.Range("B33").Value = "LidoDAOToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("E33").Value = "  +3.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.21"
$ws.Range("D33").Style = "Normal"

$ws.Range("B34").Value = "FirstDigitalUSD"
$ws.Range("C34").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("D34").Style = "Normal"

$ws.Range("E35").Value = "  -1.76%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0739"
$ws.Range("D35").Style = "Normal"

$ws.Range("E36").Value = "  +0.82%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.42"
$ws.Range("D36").Style = "Normal"

$ws.Range("E37").Value = "  -0.44%  "

$ws.Range("E38").Value = "  -2.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.105"
$ws.Range("D38").Style = "Normal"

$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("E39").Value = "  +0.37%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.82"
$ws.Range("D39").Style = "Normal"

$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("E40").Value = "  -1.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.114"
$ws.Range("D40").Style = "Normal"

$ws.Range("E41").Value = "  -4.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.01"
$ws.Range("D41").Style = "Normal"

$ws.Range("E42").Value = "  +1.38%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.30"
$ws.Range("D42").Style = "Normal"

$ws.Range("E43").Value = "  -2.56%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.965.61"
$ws.Range("D43").Style = "Normal"

$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("E44").Value = "  -4.88%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "18.92"
$ws.Range("D44").Style = "Normal"

$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("E45").Value = "  -0.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0283"
$ws.Range("D45").Style = "Normal"

$ws.Range("E46").Value = "  -2.25%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.04"
$ws.Range("D46").Style = "Normal"

$ws.Range("E47").Value = "  -1.54%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.91"
$ws.Range("D47").Style = "Normal"

$ws.Range("E48").Value = "  -0.48%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "53.45"
$ws.Range("D48").Style = "Normal"

$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("E49").Value = "  +0.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "92.37"
$ws.Range("D49").Style = "Normal"

$ws.Range("B50").Value = "BitcoinSV"
$ws.Range("C50").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("E50").Value = "  -0.70%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "71.92"
$ws.Range("D50").Style = "Normal"

$ws.Range("E51").Value = "  -1.87%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.13"
$ws.Range("D51").Style = "Normal"
